# Daily attendance processing - 2025-12-08 15:30:26
#
# For every data row, column G ("Recorded By") holds a comma-separated list
# of recorders (e.g. "System, dnasr281@gmail.com"). Re-normalize each list by
# rotating it left by one position (the first entry moves to the end), e.g.
#   "System, dnasr281@gmail.com"                  -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com, system"         -> "backup@backdoor.com, system, System"
#   "admin@admin.com, dnasr281@gmail.com"         -> "dnasr281@gmail.com, admin@admin.com"
# Single-entry lists are left untouched (rotating one element is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$recordedByCol = 7   # Column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value = $rotated
        }
    }
}
